# Apply the edit described by the diff:
#  - Insert a new "quantity" column as column B (shifting C..O right by one,
#    i.e. the old "quantity" column N becomes the new column B).
#  - Update the recalculated "Storage Fees", "Cost of Returns",
#    "Adjusted Cost of Returns", "Gross Profit" and "Gross Margin" values.
#  - Append two new summary rows: "Rich" (row 7) and "Total" (row 8).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) --------------------------------------------------
# Header cells already carry the bold/centered header style (s="1"); simply
# overwriting their text keeps that style untouched.
$ws.Cells.Item(1,1).Value = "GM"
$ws.Cells.Item(1,2).Value = "quantity"
$ws.Cells.Item(1,3).Value = "Amazon Top-line Sales (ATS)"
$ws.Cells.Item(1,4).Value = "Product Cost"
$ws.Cells.Item(1,5).Value = "Referral Fee"
$ws.Cells.Item(1,6).Value = "FBA Fulfillment Fee"
$ws.Cells.Item(1,7).Value = "other transaction fees"
$ws.Cells.Item(1,8).Value = "Shipping/Kitting Fees"
$ws.Cells.Item(1,9).Value = "Adjusted Shipping/Kitting Fees"
$ws.Cells.Item(1,10).Value = "Labeling/Polybagging Fees"
$ws.Cells.Item(1,11).Value = "Storage Fees"
$ws.Cells.Item(1,12).Value = "Allocated fees (Premium Services Fee + Subscription)"
$ws.Cells.Item(1,13).Value = "Cost of Returns"
$ws.Cells.Item(1,14).Value = "Adjusted Cost of Returns"
$ws.Cells.Item(1,15).Value = "Gross Profit"
$ws.Cells.Item(1,16).Value = "Gross Margin"

# --- Data rows (rows 2-8) -------------------------------------------------
# Force Text number format on the whole data block first so values like
# "$109,799.32", "2,666" and "19.91%" are written as literal text (matching
# the original inlineStr cells) instead of being auto-converted by Excel into
# numbers/currency/percentages.
$dataRange = $ws.Range("A2:P8")
$dataRange.NumberFormat = "@"

# Row 2: Candace
$ws.Cells.Item(2,1).Value = "Candace"
$ws.Cells.Item(2,2).Value = "2,666"
$ws.Cells.Item(2,3).Value = "`$109,799.32"
$ws.Cells.Item(2,4).Value = "`$-44,673.16"
$ws.Cells.Item(2,5).Value = "`$-16,366.10"
$ws.Cells.Item(2,6).Value = "`$-20,540.63"
$ws.Cells.Item(2,7).Value = "`$0.00"
$ws.Cells.Item(2,8).Value = "`$-2,301.96"
$ws.Cells.Item(2,9).Value = "`$-2,301.96"
$ws.Cells.Item(2,10).Value = "`$0.00"
$ws.Cells.Item(2,11).Value = "`$-3,050.59"
$ws.Cells.Item(2,12).Value = "`$-10.59"
$ws.Cells.Item(2,13).Value = "`$-1,994.42"
$ws.Cells.Item(2,14).Value = "`$-997.20"
$ws.Cells.Item(2,15).Value = "`$21,859.09"
$ws.Cells.Item(2,16).Value = "19.91%"

# Row 3: Daria
$ws.Cells.Item(3,1).Value = "Daria"
$ws.Cells.Item(3,2).Value = "190"
$ws.Cells.Item(3,3).Value = "`$2,060.64"
$ws.Cells.Item(3,4).Value = "`$-441.25"
$ws.Cells.Item(3,5).Value = "`$-309.43"
$ws.Cells.Item(3,6).Value = "`$-740.32"
$ws.Cells.Item(3,7).Value = "`$0.00"
$ws.Cells.Item(3,8).Value = "`$-18.91"
$ws.Cells.Item(3,9).Value = "`$0.00"
$ws.Cells.Item(3,10).Value = "`$-0.55"
$ws.Cells.Item(3,11).Value = "`$-175.52"
$ws.Cells.Item(3,12).Value = "`$-0.18"
$ws.Cells.Item(3,13).Value = "`$-70.71"
$ws.Cells.Item(3,14).Value = "`$-17.68"
$ws.Cells.Item(3,15).Value = "`$375.71"
$ws.Cells.Item(3,16).Value = "18.23%"

# Row 4: David M
$ws.Cells.Item(4,1).Value = "David M"
$ws.Cells.Item(4,2).Value = "2"
$ws.Cells.Item(4,3).Value = "`$47.98"
$ws.Cells.Item(4,4).Value = "`$-28.01"
$ws.Cells.Item(4,5).Value = "`$-7.20"
$ws.Cells.Item(4,6).Value = "`$-8.67"
$ws.Cells.Item(4,7).Value = "`$0.00"
$ws.Cells.Item(4,8).Value = "`$0.00"
$ws.Cells.Item(4,9).Value = "`$0.00"
$ws.Cells.Item(4,10).Value = "`$0.00"
$ws.Cells.Item(4,11).Value = "`$-0.20"
$ws.Cells.Item(4,12).Value = "`$0.00"
$ws.Cells.Item(4,13).Value = "`$0.00"
$ws.Cells.Item(4,14).Value = "`$0.00"
$ws.Cells.Item(4,15).Value = "`$3.90"
$ws.Cells.Item(4,16).Value = "8.13%"

# Row 5: David T
$ws.Cells.Item(5,1).Value = "David T"
$ws.Cells.Item(5,2).Value = "132"
$ws.Cells.Item(5,3).Value = "`$29,458.34"
$ws.Cells.Item(5,4).Value = "`$-15,049.47"
$ws.Cells.Item(5,5).Value = "`$-4,644.32"
$ws.Cells.Item(5,6).Value = "`$-821.30"
$ws.Cells.Item(5,7).Value = "`$0.00"
$ws.Cells.Item(5,8).Value = "`$-336.89"
$ws.Cells.Item(5,9).Value = "`$0.00"
$ws.Cells.Item(5,10).Value = "`$-5.61"
$ws.Cells.Item(5,11).Value = "`$-805.74"
$ws.Cells.Item(5,12).Value = "`$-2.85"
$ws.Cells.Item(5,13).Value = "`$-686.19"
$ws.Cells.Item(5,14).Value = "`$-447.55"
$ws.Cells.Item(5,15).Value = "`$7,681.50"
$ws.Cells.Item(5,16).Value = "26.08%"

# Row 6: Jocelyn
$ws.Cells.Item(6,1).Value = "Jocelyn"
$ws.Cells.Item(6,2).Value = "7,723"
$ws.Cells.Item(6,3).Value = "`$271,751.40"
$ws.Cells.Item(6,4).Value = "`$-98,771.67"
$ws.Cells.Item(6,5).Value = "`$-39,128.88"
$ws.Cells.Item(6,6).Value = "`$-55,460.44"
$ws.Cells.Item(6,7).Value = "`$0.00"
$ws.Cells.Item(6,8).Value = "`$-13,830.67"
$ws.Cells.Item(6,9).Value = "`$-13,830.67"
$ws.Cells.Item(6,10).Value = "`$-1,301.85"
$ws.Cells.Item(6,11).Value = "`$-1,346.00"
$ws.Cells.Item(6,12).Value = "`$-26.31"
$ws.Cells.Item(6,13).Value = "`$-3,925.61"
$ws.Cells.Item(6,14).Value = "`$-3,925.61"
$ws.Cells.Item(6,15).Value = "`$57,959.97"
$ws.Cells.Item(6,16).Value = "21.33%"

# Row 7: Rich
$ws.Cells.Item(7,1).Value = "Rich"
$ws.Cells.Item(7,2).Value = "0"
$ws.Cells.Item(7,3).Value = "`$0.00"
$ws.Cells.Item(7,4).Value = "`$0.00"
$ws.Cells.Item(7,5).Value = "`$0.00"
$ws.Cells.Item(7,6).Value = "`$0.00"
$ws.Cells.Item(7,7).Value = "`$0.00"
$ws.Cells.Item(7,8).Value = "`$0.00"
$ws.Cells.Item(7,9).Value = "`$0.00"
$ws.Cells.Item(7,10).Value = "`$0.00"
$ws.Cells.Item(7,11).Value = "`$-2.74"
$ws.Cells.Item(7,12).Value = "`$0.00"
$ws.Cells.Item(7,13).Value = "`$0.00"
$ws.Cells.Item(7,14).Value = "`$0.00"
$ws.Cells.Item(7,15).Value = "`$-2.74"
$ws.Cells.Item(7,16).Value = "0.00%"

# Row 8: Total
$ws.Cells.Item(8,1).Value = "Total"
$ws.Cells.Item(8,2).Value = "10,713"
$ws.Cells.Item(8,3).Value = "`$413,117.68"
$ws.Cells.Item(8,4).Value = "`$-158,963.56"
$ws.Cells.Item(8,5).Value = "`$-60,455.93"
$ws.Cells.Item(8,6).Value = "`$-77,571.36"
$ws.Cells.Item(8,7).Value = "`$0.00"
$ws.Cells.Item(8,8).Value = "`$-16,488.43"
$ws.Cells.Item(8,9).Value = "`$-16,132.63"
$ws.Cells.Item(8,10).Value = "`$-1,308.01"
$ws.Cells.Item(8,11).Value = "`$-5,380.79"
$ws.Cells.Item(8,12).Value = "`$-39.93"
$ws.Cells.Item(8,13).Value = "`$-6,676.93"
$ws.Cells.Item(8,14).Value = "`$-5,388.04"
$ws.Cells.Item(8,15).Value = "`$87,877.43"
$ws.Cells.Item(8,16).Value = "21.27%"

# Reset the data block back to the default "Normal" style (no explicit
# number-format override), matching the target workbook where data rows have
# no "s" style attribute at all.
$dataRange.Style = "Normal"

Write-Output "done"
